$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old standalone "Docentes responsaveis" data row (row 13: B13/C13
# held the professor name with no label in column A). Deleting it shifts all
# subsequent rows up by one, matching the target layout (A1:C22 -> A1:C21).
$ws.Rows.Item(13).Delete()

$ws.Range("B10").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C10").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2015"
$ws.Range("C15").Value = "01/01/2015"

$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("B19").Value = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.`n`nOs alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. `nCada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.`nAs aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas."
$ws.Range("C19").Value = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.`n`nOs alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. `nCada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.`nAs aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas."

$ws.Range("B20").Value = "A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.`nO detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina."
$ws.Range("C20").Value = "A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.`nO detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina."

$ws.Range("B21").Value = "Não há recuperação"
$ws.Range("C21").Value = "Não há recuperação"
